$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E; existing D:K shifts to F:M
$ws.Columns("D:E").Insert()

# Copy number formats from column F (the old column D, now shifted) into new D:E columns
$ws.Range("F5:F102").Copy() | Out-Null
$ws.Range("D5:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the two new quarter columns (D = most recent quarter, E = previous quarter)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1297000
$ws.Range("E8").Value = 1239000
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = -22000
$ws.Range("E15").Value = -23000
$ws.Range("D17").Value = 356000
$ws.Range("E17").Value = 315000
$ws.Range("D18").Value = 941000
$ws.Range("E18").Value = 924000
$ws.Range("D20").Value = -367000
$ws.Range("E20").Value = -355000
$ws.Range("D21").Value = 666000
$ws.Range("E21").Value = 664000
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 574000
$ws.Range("E23").Value = 569000
$ws.Range("D24").Value = 92000
$ws.Range("E24").Value = 87000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 482000
$ws.Range("E26").Value = 482000
$ws.Range("D27").Value = 459000
$ws.Range("E27").Value = 468000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 2000
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 367000
$ws.Range("E32").Value = 355000
$ws.Range("D33").Value = 461000
$ws.Range("E33").Value = 468000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 461000
$ws.Range("E35").Value = 468000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 678000
$ws.Range("E41").Value = 319000
$ws.Range("D42").Value = 1515000
$ws.Range("E42").Value = 1639000
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 882000
$ws.Range("E48").Value = 891000
$ws.Range("D49").Value = 2832000
$ws.Range("E49").Value = 2854000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1100000
$ws.Range("E52").Value = 1488000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 139613000
$ws.Range("E54").Value = 138805000
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 1421000
$ws.Range("E59").Value = 2044000
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 13732000
$ws.Range("E61").Value = 13849000
$ws.Range("D62").Value = 692000
$ws.Range("E62").Value = "NA"
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 124018000
$ws.Range("E66").Value = 123597000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 1450000
$ws.Range("E70").Value = 1450000
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 11556000
$ws.Range("E72").Value = 11262000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 14145000
$ws.Range("E76").Value = 13758000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 461000
$ws.Range("E81").Value = 468000
$ws.Range("D83").Value = 92000
$ws.Range("E83").Value = 95000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 1712000
$ws.Range("E89").Value = 386000
$ws.Range("D91").Value = -26000
$ws.Range("E91").Value = -37000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -1099000
$ws.Range("E94").Value = -1424000
$ws.Range("D96").Value = -198000
$ws.Range("E96").Value = -190000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -254000
$ws.Range("E100").Value = 573000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 359000
$ws.Range("E102").Value = -465000
